# Generate Report for Handoff
# Updates the localization-status report to reflect that b.md is now
# "Ready for handoff": refreshes status/handoff metadata for the b.md
# row on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f410ceb6e1d36ec87a166dc3987e15bb7eb0daf0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef48d997754203fd630c3ca423817f104b95f25f/e2e/b.md."

# --- Overview sheet: b.md row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-20 08:46:25"

# --- zh-cn sheet: b.md row (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# "Content Duplicate" is stored as literal text ("True"/"False"), not a
# boolean, elsewhere in this report. Copy/PasteSpecial-values from the
# existing text "False" cell (F2) so F3 keeps the same text representation
# instead of being auto-coerced into a native boolean by a plain .Value=.
$wsZhCn.Range("F2").Copy()
$wsZhCn.Range("F3").PasteSpecial(-4163)  # xlPasteValues
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-20 08:46:20"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet: b.md row (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F2").Copy()
$wsDeDe.Range("F3").PasteSpecial(-4163)  # xlPasteValues
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-20 08:46:25"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 40
